$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.2149532710280374
$ws.Cells.Item(2, 3).Value = 0.4953271028037383
$ws.Cells.Item(2, 10).Value = 0.01869158878504673
$ws.Cells.Item(2, 16).Value = 0.1448598130841121
$ws.Cells.Item(2, 19).Value = 0.1261682242990654
$ws.Cells.Item(3, 2).Value = 0.009345794392523364
$ws.Cells.Item(3, 10).Value = 0.04672897196261682
$ws.Cells.Item(3, 16).Value = 0.6635514018691588
$ws.Cells.Item(3, 19).Value = 0.2803738317757009
$ws.Cells.Item(4, 10).Value = 0.02325581395348837
$ws.Cells.Item(4, 16).Value = 0.6976744186046512
$ws.Cells.Item(4, 19).Value = 0.2790697674418605
$ws.Cells.Item(6, 2).Value = 0.03864734299516908
$ws.Cells.Item(6, 4).Value = 0.01932367149758454
$ws.Cells.Item(6, 6).Value = 0.03381642512077294
$ws.Cells.Item(6, 10).Value = 0.3961352657004831
$ws.Cells.Item(6, 15).Value = 0.00966183574879227
$ws.Cells.Item(6, 17).Value = 0.1835748792270532
$ws.Cells.Item(6, 18).Value = 0.04830917874396135
$ws.Cells.Item(6, 19).Value = 0.2705314009661836
$ws.Cells.Item(7, 2).Value = 0.1258278145695364
$ws.Cells.Item(7, 4).Value = 0.01324503311258278
$ws.Cells.Item(7, 6).Value = 0.05960264900662252
$ws.Cells.Item(7, 10).Value = 0.1059602649006623
$ws.Cells.Item(7, 15).Value = 0.006622516556291391
$ws.Cells.Item(7, 17).Value = 0.1788079470198675
$ws.Cells.Item(7, 18).Value = 0.1059602649006623
$ws.Cells.Item(7, 19).Value = 0.4039735099337748
$ws.Cells.Item(8, 2).Value = 0.05882352941176471
$ws.Cells.Item(8, 4).Value = 0.02450980392156863
$ws.Cells.Item(8, 6).Value = 0.06127450980392157
$ws.Cells.Item(8, 10).Value = 0.1151960784313725
$ws.Cells.Item(8, 15).Value = 0.01225490196078431
$ws.Cells.Item(8, 17).Value = 0.1764705882352941
$ws.Cells.Item(8, 18).Value = 0.1519607843137255
$ws.Cells.Item(8, 19).Value = 0.3995098039215687
$ws.Cells.Item(9, 2).Value = 0.05797101449275362
$ws.Cells.Item(9, 4).Value = 0.01449275362318841
$ws.Cells.Item(9, 6).Value = 0.07246376811594203
$ws.Cells.Item(9, 10).Value = 0.09420289855072464
$ws.Cells.Item(9, 15).Value = 0.02173913043478261
$ws.Cells.Item(9, 17).Value = 0.2028985507246377
$ws.Cells.Item(9, 18).Value = 0.1594202898550725
$ws.Cells.Item(9, 19).Value = 0.3768115942028986
$ws.Cells.Item(10, 2).Value = 0.09228039041703638
$ws.Cells.Item(10, 4).Value = 0.02218278615794144
$ws.Cells.Item(10, 5).Value = 0.0008873114463176575
$ws.Cells.Item(10, 6).Value = 0.06122448979591837
$ws.Cells.Item(10, 10).Value = 0.1118012422360248
$ws.Cells.Item(10, 15).Value = 0.01508429458740018
$ws.Cells.Item(10, 17).Value = 0.1925465838509317
$ws.Cells.Item(10, 18).Value = 0.1224489795918367
$ws.Cells.Item(10, 19).Value = 0.3815439219165927
$ws.Cells.Item(11, 7).Value = 0.1416666666666667
$ws.Cells.Item(11, 10).Value = 0.1083333333333333
$ws.Cells.Item(11, 11).Value = 0.1958333333333333
$ws.Cells.Item(11, 12).Value = 0.55
$ws.Cells.Item(11, 19).Value = 0.004166666666666667
$ws.Cells.Item(12, 7).Value = 0.7348484848484849
$ws.Cells.Item(12, 10).Value = 0.1893939393939394
$ws.Cells.Item(12, 11).Value = 0.007575757575757576
$ws.Cells.Item(12, 12).Value = 0.03787878787878788
$ws.Cells.Item(12, 19).Value = 0.0303030303030303
$ws.Cells.Item(13, 7).Value = 0.7058823529411765
$ws.Cells.Item(13, 10).Value = 0.2941176470588235
$ws.Cells.Item(15, 6).Value = 0.03076923076923077
$ws.Cells.Item(15, 8).Value = 0.1794871794871795
$ws.Cells.Item(15, 9).Value = 0.04102564102564103
$ws.Cells.Item(15, 10).Value = 0.4153846153846154
$ws.Cells.Item(15, 11).Value = 0.06153846153846154
$ws.Cells.Item(15, 13).Value = 0.03076923076923077
$ws.Cells.Item(15, 15).Value = 0.03589743589743589
$ws.Cells.Item(15, 19).Value = 0.2051282051282051
$ws.Cells.Item(16, 6).Value = 0.03875968992248062
$ws.Cells.Item(16, 8).Value = 0.2403100775193799
$ws.Cells.Item(16, 9).Value = 0.08527131782945736
$ws.Cells.Item(16, 10).Value = 0.3178294573643411
$ws.Cells.Item(16, 11).Value = 0.1317829457364341
$ws.Cells.Item(16, 13).Value = 0.02325581395348837
$ws.Cells.Item(16, 15).Value = 0.02325581395348837
$ws.Cells.Item(16, 19).Value = 0.1395348837209302
$ws.Cells.Item(17, 6).Value = 0.03149606299212598
$ws.Cells.Item(17, 8).Value = 0.1916010498687664
$ws.Cells.Item(17, 9).Value = 0.05774278215223097
$ws.Cells.Item(17, 10).Value = 0.4514435695538058
$ws.Cells.Item(17, 11).Value = 0.1181102362204724
$ws.Cells.Item(17, 13).Value = 0.005249343832020997
$ws.Cells.Item(17, 15).Value = 0.06561679790026247
$ws.Cells.Item(17, 19).Value = 0.07874015748031496
$ws.Cells.Item(18, 6).Value = 0.0564516129032258
$ws.Cells.Item(18, 8).Value = 0.157258064516129
$ws.Cells.Item(18, 9).Value = 0.04838709677419355
$ws.Cells.Item(18, 10).Value = 0.4879032258064516
$ws.Cells.Item(18, 11).Value = 0.09677419354838709
$ws.Cells.Item(18, 13).Value = 0.01209677419354839
$ws.Cells.Item(18, 15).Value = 0.06451612903225806
$ws.Cells.Item(18, 19).Value = 0.07661290322580645
$ws.Cells.Item(19, 6).Value = 0.03288201160541586
$ws.Cells.Item(19, 8).Value = 0.2263056092843327
$ws.Cells.Item(19, 9).Value = 0.08220502901353965
$ws.Cells.Item(19, 10).Value = 0.3646034816247582
$ws.Cells.Item(19, 11).Value = 0.0851063829787234
$ws.Cells.Item(19, 13).Value = 0.02030947775628627
$ws.Cells.Item(19, 14).Value = 0.001934235976789168
$ws.Cells.Item(19, 15).Value = 0.08123791102514506
$ws.Cells.Item(19, 19).Value = 0.1054158607350097
